# 3.1 Bond Valuation Example.xlsx
#
# The "Template" worksheet (the blank exercise students fill in) is
# duplicated into a new "Answer" worksheet (the worked solution, with the
# final bond price highlighted in yellow), and the Template is then reset
# to a fresh "today" / first coupon date so the exercise no longer matches
# the answer key verbatim.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Template")

# --- 1. Duplicate the Template sheet (this captures the *current* state of
#        Template - old "today" anchor + old first-coupon date - and that
#        copy becomes the "Answer" sheet) -------------------------------
$template.Copy([System.Reflection.Missing]::Value, $template)
$answer = $wb.Worksheets.Item($template.Index + 1)
$answer.Name = "Answer"

# The answer key has no "Bond Price" caption next to its total - only the
# total itself, highlighted in yellow.
$answer.Range("D10").Clear()
$answer.Range("E10").Interior.Color = 65535   # RGB(255,255,0) -> yellow

# --- 2. Reset the Template sheet back to a blank-ish exercise -----------
# "Today" now reads two days earlier than a fresh TODAY() would, so the
# template keeps showing the original reference date.
$template.Range("B1").Formula = "=TODAY()-2"

# First coupon date moves from 2021-12-31 to 2021-12-15; every other date
# / discount-factor / PV cell in the sheet is a formula relative to this
# one (and to B1), so they recompute on their own.
$template.Range("A5").Value = 44545

# Put the cursor/selection back on the (now-edited) first coupon date and
# make sure Template is the tab that is on top when the file is saved.
$template.Activate()
$template.Range("A5").Select()

Write-Output "Added Answer sheet and refreshed Template dates"
